# Corrected some selection scopes
# Clear out the erroneous G/H values that were written into the odd
# data rows (3,5,7,9,11,13,15,17) of the qoq error series sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3, 5, 7, 9, 11, 13, 15, 17)

foreach ($r in $rows) {
    $ws.Range("G$r`:H$r").ClearContents()
}
